# Sentinel pricing tiers and recommendation thresholds — May 26th 2023 price/threshold refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (selectLockedCells only); lift protection so the
# price / formula cells can be edited, matching what the sheet author did
# before touching the numbers.
$ws.Unprotect()

# ---------------------------------------------------------------------
# Table 1 — "Microsoft Sentinel pricing" (rows 5-14)
# ---------------------------------------------------------------------
$ws.Range("C5").Value2 = 2.708

$ws.Range("C7").Value2 = 228.95
$ws.Range("C8").Value2 = 429.87
$ws.Range("C9").Value2 = 630.78
$ws.Range("C10").Value2 = 822.35
$ws.Range("C11").Value2 = 1010.42
$ws.Range("C12").Value2 = 1985.79
$ws.Range("C13").Value2 = 3878.12
$ws.Range("C14").Value2 = 9403.27

# Recommendation-threshold formulas no longer reference the (now cleared)
# helper column L — they recompute the effective per-GB rate inline.
$ws.Range("H8").Formula  = "=C8/(C7/D7)"
$ws.Range("H9").Formula  = "=C9/(C8/D8)"
$ws.Range("H10").Formula = "=C10/(C9/D9)"
$ws.Range("H11").Formula = "=C11/(C10/D10)"
$ws.Range("H12").Formula = "=C12/(C11/D11)"
$ws.Range("H13").Formula = "=C13/(C12/D12)"
$ws.Range("H14").Formula = "=C14/(C13/D13)"

# H14 previously carried its own right-aligned number style; normalise it to
# match the rest of the column (same as H29 below).
$ws.Range("H14").HorizontalAlignment = 1

# ---------------------------------------------------------------------
# Table 2 — "Log Analytics / Azure Monitor" (rows 20-29)
# ---------------------------------------------------------------------
$ws.Range("C20").Value2 = 2.36

$ws.Range("C22").Value2 = 117.72
$ws.Range("C23").Value2 = 211.89
$ws.Range("C24").Value2 = 306.07
$ws.Range("C25").Value2 = 392.39
$ws.Range("C26").Value2 = 470.87
$ws.Range("C27").Value2 = 918.19
$ws.Range("C28").Value2 = 1742.21
$ws.Range("C29").Value2 = 4120.08

$ws.Range("H23").Formula = "=C23/(C22/D22)"
$ws.Range("H24").Formula = "=C24/(C23/D23)"
$ws.Range("H25").Formula = "=C25/(C24/D24)"
$ws.Range("H26").Formula = "=C26/(C25/D25)"
$ws.Range("H27").Formula = "=C27/(C26/D26)"
$ws.Range("H28").Formula = "=C28/(C27/D27)"
$ws.Range("H29").Formula = "=C29/(C28/D28)"

# ---------------------------------------------------------------------
# Retire the helper "Price/GB" column (L): clear its header + the per-tier
# helper values now that H uses inline ratios, and unhide the column.
# ---------------------------------------------------------------------
$ws.Range("L5").ClearContents()
$ws.Range("L7:L14").ClearContents()
$ws.Range("L22:L29").ClearContents()
$ws.Columns("L").Hidden = $false

# ---------------------------------------------------------------------
# Small formatting leftovers picked up while touching the tables (copy the
# format from an already-styled neighbour so the underlying number format
# matches exactly).
# ---------------------------------------------------------------------
$ws.Range("R8").NumberFormat  = $ws.Range("J22").NumberFormat
$ws.Range("Q29").NumberFormat = $ws.Range("L34").NumberFormat

# Final selection left on the last recomputed threshold cell.
$ws.Range("H14").Select()
